$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.156242251396179
$ws.Range("B1").Value = 2.064386129379272
$ws.Range("C1").Value = 3.139245986938477
$ws.Range("D1").Value = 1.037437081336975
$ws.Range("E1").Value = 1.51601231098175
